$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.295.76"
$ws.Range("E2").Value = "  +0.01%  "

$ws.Range("D3").Value = "2.963.86"
$ws.Range("E3").Value = "  +2.06%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "357.30"
$ws.Range("E5").Value = "  +1.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.53"
$ws.Range("E6").Value = "  -3.64%  "

$ws.Range("E7").Value = "  +2.22%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("E9").Value = "  +0.51%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.01"
$ws.Range("E10").Value = "  -2.72%  "

$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0878"
$ws.Range("E11").Value = "  +1.35%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.138"
$ws.Range("E12").Value = "  +1.52%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.57"
$ws.Range("E13").Value = "  -1.33%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "3.409.56"
$ws.Range("E14").Value = "  +1.26%  "

$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.77"
$ws.Range("E15").Value = "  +0.05%  "

$ws.Range("D16").Value = "2.955.07"
$ws.Range("E16").Value = "  +1.22%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.984"
$ws.Range("E17").Value = "  -1.45%  "

$ws.Range("D18").Value = "52.286.07"
$ws.Range("E18").Value = "  -0.09%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.51"
$ws.Range("E19").Value = "  +6.42%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.67"
$ws.Range("E20").Value = "  +0.17%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.88"
$ws.Range("E21").Value = "  -1.91%  "

$ws.Range("E22").Value = "  +0.62%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.40"
$ws.Range("E23").Value = "  -0.61%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "270.86"
$ws.Range("E24").Value = "  +0.61%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.80"
$ws.Range("E25").Value = "  +0.35%  "

$ws.Range("E26").Value = "  +2.64%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.83"
$ws.Range("E27").Value = "  +19.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "27.08"
$ws.Range("E28").Value = "  +1.06%  "

$ws.Range("E29").Value = "  +0.06%  "

$ws.Range("E30").Value = "  +2.20%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "10.50"
$ws.Range("E31").Value = "  -1.12%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "37.79"
$ws.Range("E32").Value = "  +0.31%  "

$ws.Range("B33").Value = "RenderToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.19"
$ws.Range("E33").Value = "  -1.16%  "

$ws.Range("B34").Value = "Toncoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.19"
$ws.Range("E34").Value = "  +11.00%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "52.11"
$ws.Range("E35").Value = "  -2.25%  "

$ws.Range("E36").Value = "  -1.63%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.05%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.22"
$ws.Range("E38").Value = "  -3.56%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.22"
$ws.Range("E39").Value = "  -3.67%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.00"
$ws.Range("E40").Value = "  -2.98%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.72"
$ws.Range("E41").Value = "  -1.51%  "

$ws.Range("E42").Value = "  +2.87%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "23.32"
$ws.Range("E43").Value = "  +0.68%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "119.21"
$ws.Range("E44").Value = "  -0.73%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.18"
$ws.Range("E45").Value = "  -0.86%  "

$ws.Range("E46").Value = "  -5.19%  "

$ws.Range("E47").Value = "  -2.00%  "

$ws.Range("D48").Value = "2.143.87"
$ws.Range("E48").Value = "  -1.81%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.248"
$ws.Range("E49").Value = "  -5.63%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0354"
$ws.Range("E50").Value = "  +1.13%  "

$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.37"
$ws.Range("E51").Value = "  +1.49%  "
